# Limited the number of Skills and traits
# Also updates student name / target-institution fields on the Student Profile sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Student identity fields -------------------------------------------------
$ws.Range("B3").Value = "Trei"
$ws.Range("B4").Value = "Tulia"
$ws.Range("B9").Value = "MIT"

# --- Remove stray formatted-but-empty cells (B10, B18, and all of row 24) ----
$ws.Range("B10").Clear()
$ws.Range("B18").Clear()
$ws.Range("A24:B24").Clear()

# Those two rows (10 & 18) had borders applied to column A that are no longer
# needed once the neighboring empty cell is gone - clear them to match.
$ws.Range("A10").Borders.LineStyle = 0
$ws.Range("A18").Borders.LineStyle = 0

# --- Purpose of letter: "University" selected instead of "Other" ------------
$ws.Range("B14").Value = "X"
$ws.Range("B16").ClearContents()

# --- Positive Personality Traits: limit selection down to 3 -----------------
# Previously checked: altruistic, assertive, amiable, brilliant, bright, determined
# Now checked:        assertive, enthusiastic, logical
$ws.Range("B27").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("B35").ClearContents()
$ws.Range("B36").Value = "X"
$ws.Range("B44").Value = "X"

# --- Academic Skills: limit selection down to 3 ------------------------------
# Previously checked: disciplined work habits, problem-solving skills, teamwork,
#                      adaptability to new environment, creative, great presentation skills
# Now checked:         disciplined work habits, leadership skills, great presentation skills
$ws.Range("B70").ClearContents()
$ws.Range("B72").ClearContents()
$ws.Range("B73").Value = "X"
$ws.Range("B74").ClearContents()
$ws.Range("B76").ClearContents()

# --- Update selection to match the author's last cursor position ------------
$ws.Range("B15").Select()
